$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Insert a new row above row 13 (shifts existing rows 13-23 down to 14-24)
$ws.Rows.Item(13).Insert()

# New row 13 data
$ws.Cells.Item(13, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Cells.Item(13, 2).Value = 11128
$ws.Cells.Item(13, 4).Value = 100
$ws.Cells.Item(13, 5).Value = $null

# Updated counts for rows that shifted down (14-18), values/labels/dates already
# moved with the row insert, only the Client/Sample Count (col B) needs updating
$ws.Cells.Item(14, 2).Value = 486214
$ws.Cells.Item(15, 2).Value = 79953
$ws.Cells.Item(16, 2).Value = 35355
$ws.Cells.Item(17, 2).Value = 65425
$ws.Cells.Item(18, 2).Value = 117653
